$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ01319359",
    "summ01443819",
    "summ01569287",
    "summ01681272",
    "summ01790130",
    "summ01902538",
    "summ02014848",
    "summ02118752",
    "summ02232015",
    "summ02343920",
    "summ02452467",
    "summ02568254",
    "summ02675354",
    "summ02789293",
    "summ02900874",
    "summ03011783",
    "summ03128970",
    "summ03261968",
    "summ03415455",
    "summ03574817",
    "summ03725965",
    "summ03920943",
    "summ04098836",
    "summ04278905",
    "summ04453445",
    "summ04641131",
    "summ04829886",
    "summ05004498",
    "summ05154973",
    "summ05309897",
    "summ05457771",
    "summ05604116",
    "summ05786958",
    "summ05968031",
    "summ06129371",
    "summ06279101",
    "summ06441137",
    "summ06615117",
    "summ06776367",
    "summ06927420",
    "summ07133643",
    "summ07292148",
    "summ07449600",
    "summ07601347",
    "summ07763778",
    "summ07945161",
    "summ08108533",
    "summ08258140",
    "summ08421051",
    "summ08578059"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
